$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells are treated as text so values such as
# "248.58" are not auto-converted to numbers by Excel (matches source data
# which stores these as plain strings).
$ws.Range("D2:D51").NumberFormat = "@"

# Updated price/volume figures for existing rows
$ws.Range("D2").Value = '37.109.12'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '2.055.85'
$ws.Range("E3").Value = '  -2.68%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '248.58'
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '55.44'
$ws.Range("E8").Value = '  +17.23%  '
$ws.Range("D9").Value = '61.44'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").Value = '0.0790'
$ws.Range("E12").Value = '  +5.53%  '
$ws.Range("E13").Value = '  +6.41%  '
$ws.Range("D14").Value = '2.352.53'
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = '0.818'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").Value = '2.056.10'
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D18").Value = '37.059.03'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").Value = '72.36'
$ws.Range("E19").Value = '  -1.64%  '
$ws.Range("D20").Value = '0.0₃0902'
$ws.Range("E20").Value = '  +8.00%  '
$ws.Range("D21").Value = '14.18'
$ws.Range("E21").Value = '  +7.33%  '
$ws.Range("D22").Value = '5.32'
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").Value = '236.84'
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("D26").Value = '169.77'
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("D28").Value = '20.12'
$ws.Range("E28").Value = '  -6.81%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("E32").Value = '  +10.40%  '
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("E34").Value = '  +4.59%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = '0.0861'
$ws.Range("E36").Value = '  -8.59%  '
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  -6.08%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  +22.19%  '
$ws.Range("D41").Value = '18.13'
$ws.Range("E41").Value = '  +13.69%  '
$ws.Range("E42").Value = '  -0.64%  '
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("D44").Value = '95.83'
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D48").Value = '2.42'
$ws.Range("E48").Value = '  +5.96%  '
$ws.Range("D49").Value = '1.295.76'
$ws.Range("E49").Value = '  -3.97%  '
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("D51").Value = '6.78'
$ws.Range("E51").Value = '  -5.61%  '

# Rows 45-47 reshuffled: Gas moves up to rank 45, FTXToken to 46, HuobiToken to 47
$ws.Range("B45").Value = 'Gas'
$ws.Range("C45").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D45").Value = '14.73'
$ws.Range("E45").Value = '  -49.40%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.23'
$ws.Range("E46").Value = '  +46.01%  '

$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").Value = '2.77'
$ws.Range("E47").Value = '  -0.44%  '
